$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.722.38"
$ws.Range("E2").Value = "  -0.30%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.295.33"
$ws.Range("E3").Value = "  -0.05%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "305.15"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "96.30"
$ws.Range("E6").Value = "  -0.98%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.504"
$ws.Range("E7").Value = "  -1.74%  "
$ws.Range("E8").Value = "  +0.08%  "
$ws.Range("E9").Value = "  -1.68%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.04"
$ws.Range("E10").Value = "  -2.68%  "
$ws.Range("E11").Value = "  -0.59%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "18.63"
$ws.Range("E12").Value = "  +5.07%  "
$ws.Range("E13").Value = "  +2.17%  "
$ws.Range("E14").Value = "  +1.26%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.653.09"
$ws.Range("E15").Value = "  +0.24%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.322.92"
$ws.Range("E16").Value = "  +0.95%  "
$ws.Range("E17").Value = "  -0.06%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "42.657.73"
$ws.Range("E18").Value = "  -0.32%  "
$ws.Range("E19").Value = "  +1.45%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0892"
$ws.Range("E20").Value = "  -1.64%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.01"
$ws.Range("E21").Value = "  -0.87%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "67.11"
$ws.Range("E22").Value = "  -1.16%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "235.51"
$ws.Range("E23").Value = "  -2.64%  "
$ws.Range("E24").Value = "  +0.77%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.41"
$ws.Range("E26").Value = "  -0.46%  "
$ws.Range("E27").Value = "  -2.16%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "166.35"
$ws.Range("E28").Value = "  +0.21%  "
$ws.Range("E29").Value = "  +0.98%  "
$ws.Range("E30").Value = "  -0.03%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "32.82"
$ws.Range("E31").Value = "  +0.05%  "
$ws.Range("E32").Value = "  +0.07%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "18.00"
$ws.Range("E33").Value = "  +4.91%  "
$ws.Range("E34").Value = "  -0.52%  "
$ws.Range("B35").Value = "RenderToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.44"
$ws.Range("E35").Value = "  -7.13%  "
$ws.Range("B36").Value = "WEMIXToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.35"
$ws.Range("E36").Value = "  -1.22%  "
$ws.Range("E37").Value = "  -0.26%  "
$ws.Range("E39").Value = "  -1.15%  "
$ws.Range("E40").Value = "  -0.81%  "
$ws.Range("E41").Value = "  -1.76%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.997.24"
$ws.Range("E42").Value = "  -0.71%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0280"
$ws.Range("E43").Value = "  -1.36%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "10.25"
$ws.Range("E44").Value = "  +1.30%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "18.08"
$ws.Range("E45").Value = "  +5.20%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.11"
$ws.Range("E46").Value = "  -0.37%  "
$ws.Range("E47").Value = "  -0.45%  "
$ws.Range("B48").Value = "MultiversX"
$ws.Range("C48").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "53.41"
$ws.Range("E48").Value = "  +0.50%  "
$ws.Range("B49").Value = "HuobiToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.83"
$ws.Range("E49").Value = "  +2.22%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.520.55"
$ws.Range("E50").Value = "  -0.15%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "70.83"
$ws.Range("E51").Value = "  -1.62%  "
